$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 1.62
$ws.Range("H2").Value = 3.9
$ws.Range("I2").Value = 5.5
$ws.Range("J2").Value = 2.2
$ws.Range("L2").Value = 5.5
$ws.Range("U2").Value = 1.91
$ws.Range("V2").Value = 1.91
$ws.Range("X2").Value = 7.5
$ws.Range("Z2").Value = 12
$ws.Range("AG2").Value = 301
$ws.Range("AO2").Value = 8
$ws.Range("AS2").Value = 151
$ws.Range("AU2").Value = 8.5
$ws.Range("AW2").Value = 7
$ws.Range("AX2").Value = 29
$ws.Range("AZ2").Value = 101
$ws.Range("BA2").Value = 126
$ws.Range("BB2").Value = 251
$ws.Range("M4").Value = 1.13
$ws.Range("N4").Value = 6
$ws.Range("X4").Value = 6.5
$ws.Range("AA4").Value = 21
$ws.Range("AS4").Value = 351
$ws.Range("AU4").Value = 11
$ws.Range("AV4").Value = 101
$ws.Range("Q5").Value = 1.98
$ws.Range("R5").Value = 1.88
$ws.Range("O6").Value = 1.11
$ws.Range("P6").Value = 4.6
$ws.Range("S6").Value = 1.26
$ws.Range("T6").Value = 3.48
$ws.Range("U6").Value = 1.81
$ws.Range("V6").Value = 1.95
$ws.Range("H7").Value = 2.8
$ws.Range("I7").Value = 2.5
$ws.Range("J7").Value = 3.5
$ws.Range("N7").Value = 6.95
$ws.Range("O7").Value = 1.37
$ws.Range("P7").Value = 2.62
$ws.Range("Q7").Value = 2.1
$ws.Range("R7").Value = 1.57
$ws.Range("S7").Value = 1.42
$ws.Range("T7").Value = 2.47
$ws.Range("U7").Value = 1.75
$ws.Range("V7").Value = 1.85
$ws.Range("AC7").Value = 7.3
$ws.Range("AD7").Value = 5.5
$ws.Range("AH7").Value = 7.2
$ws.Range("AJ7").Value = 9.25
$ws.Range("AL7").Value = 23
$ws.Range("AM7").Value = 32
$ws.Range("AO7").Value = 16.5
$ws.Range("AP7").Value = 22
$ws.Range("AT7").Value = 2.45
$ws.Range("AU7").Value = 6.4
$ws.Range("AV7").Value = 55
$ws.Range("AZ7").Value = 55
$ws.Range("G9").Value = 1.45
$ws.Range("H9").Value = 4.5
$ws.Range("I9").Value = 6.5
$ws.Range("J9").Value = 1.95
$ws.Range("K9").Value = 2.6
$ws.Range("Q9").Value = 1.5
$ws.Range("R9").Value = 2.5
$ws.Range("AD9").Value = 9
$ws.Range("AQ9").Value = 19
$ws.Range("AW9").Value = 8
$ws.Range("AZ9").Value = 101
$ws.Range("BD9").Value = 151
$ws.Range("Q10").Value = 2.3
$ws.Range("R10").Value = 1.6
$ws.Range("M11").Value = 1.11
$ws.Range("N11").Value = 6.5
$ws.Range("G13").Value = 1.29
$ws.Range("H13").Value = 5
$ws.Range("I13").Value = 12
$ws.Range("O13").Value = 1.2
$ws.Range("P13").Value = 4.33
$ws.Range("Q13").Value = 1.67
$ws.Range("R13").Value = 2.15
$ws.Range("U13").Value = 2.2
$ws.Range("V13").Value = 1.62
$ws.Range("Y13").Value = 9.5
$ws.Range("AA13").Value = 12
$ws.Range("AC13").Value = 12
$ws.Range("AH13").Value = 23
$ws.Range("AJ13").Value = 29
$ws.Range("AS13").Value = 151
$ws.Range("AZ13").Value = 251
$ws.Range("BA13").Value = 251
